$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-recognized as numbers by Excel, so they remain text like the originals.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D13", "D14", "D15", "D16", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D49")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated values scraped for this run.
$ws.Range("D2").Value = "41.926.72"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "2.208.88"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "230.27"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "60.67"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "0.0901"
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "2.537.58"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "15.42"
$ws.Range("E13").Value = "  -2.94%  "
$ws.Range("D14").Value = "22.04"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "0.798"
$ws.Range("D16").Value = "5.57"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "2.210.13"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "41.865.34"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "0.0₃0943"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("D20").Value = "72.20"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("D21").Value = "6.06"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").Value = "242.53"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("D26").Value = "9.62"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "168.87"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("D29").Value = "20.30"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  -1.54%  "
$ws.Range("E31").Value = "  -5.11%  "
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -4.90%  "
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("D35").Value = "0.0646"
$ws.Range("E35").Value = "  +3.02%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "3.54"
$ws.Range("E36").Value = "  -7.57%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "6.27"
$ws.Range("E37").Value = "  -6.87%  "
$ws.Range("D38").Value = "2.33"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0244"
$ws.Range("E39").Value = "  +2.61%  "
$ws.Range("B40").Value = "BinanceUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "0.000225"
$ws.Range("E41").Value = "  -10.61%  "
$ws.Range("D42").Value = "8.54"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").Value = "0.0952"
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "96.55"
$ws.Range("E45").Value = "  -4.06%  "
$ws.Range("B46").Value = "FTXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D46").Value = "4.37"
$ws.Range("E46").Value = "  -12.62%  "
$ws.Range("D47").Value = "1.454.91"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("D49").Value = "16.06"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("E51").Value = "  +1.55%  "
